$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.339.64"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.910.81"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +8.89%  "
$ws.Range("D6").Value = "'254.17"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'40.81"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "'0.359"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "'0.0750"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "2.189.97"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "'12.64"
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("E15").Value = "  +3.23%  "
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "1.912.91"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "35.360.51"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'74.32"
$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").Value = "'244.11"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("D23").Value = "'5.06"
$ws.Range("E23").Value = "  +4.70%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'2.46"
$ws.Range("E25").Value = "  +5.82%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").Value = "'166.84"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "'8.64"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").Value = "'0.133"
$ws.Range("E30").Value = "  +4.64%  "
$ws.Range("D31").Value = "4.130.28"
$ws.Range("E31").Value = "  +19.50%  "
$ws.Range("D32").Value = "'4.35"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("E33").Value = "  +14.82%  "
$ws.Range("E34").Value = "  +22.55%  "
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'0.912"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.28"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0219"
$ws.Range("E41").Value = "  +4.79%  "
$ws.Range("D42").Value = "'97.16"
$ws.Range("E42").Value = "  +7.67%  "
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.340.77"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'45.24"
$ws.Range("E50").Value = "  -6.43%  "
$ws.Range("D51").Value = "'12.12"
$ws.Range("E51").Value = "  +17.64%  "
